{"js": "// Fix a typo in the SmartArt (\"Archivist\") diagram embedded in the\n// document:\n//   \"Aufsetzen des Systems auf projektfremden einem Rechner\"\n//     -> \"Aufsetzen des Systems auf projektfremdem Rechner\"\n\nconst oldFragment = \"projektfremden einem Rechner\";\nconst newFragment = \"projektfremdem Rechner\";\nconst oldText = \"Aufsetzen des Systems auf projektfremden einem Rechner\";\nconst newText = \"Aufsetzen des Systems auf projektfremdem Rechner\";\n\nlet fixed = false;\n\n// 1) The wording lives inside a SmartArt graphic. Office.js has no\n//    dedicated SmartArt text API, but shapes (incl. SmartArt graphic\n//    frames) are enumerable off the body - walk them and patch any\n//    exposed text surface we can find.\nconst shapes = context.document.body.shapes;\nshapes.load(\"items\");\nawait context.sync();\n\nfor (const shape of shapes.items) {\n  // Best-effort: some hosts expose a text body / text-frame style API on\n  // a shape; only touch it if present so this stays a no-op elsewhere.\n  const anyShape = shape;\n  try {\n    if (typeof anyShape.getTextFrame === \"function\") {\n      const tf = anyShape.getTextFrame();\n      tf.load(\"text\");\n      await context.sync();\n      const t = tf.text;\n      if (t && t.indexOf(oldFragment) !== -1) {\n        tf.deleteText ? tf.deleteText() : null;\n        if (typeof tf.insertText === \"function\") {\n          tf.insertText(t.split(oldFragment).join(newFragment), \"Replace\");\n          await context.sync();\n          fixed = true;\n        }\n      }\n    }\n  } catch (e) {\n    // Shape has no text surface in this host - ignore and move on.\n  }\n}\n\n// 2) Belt-and-braces: run an ordinary body-wide search & replace too, so\n//    the fix still lands if the phrase ever lives in normal paragraph\n//    text instead of (or in addition to) the diagram.\nconst body = context.document.body;\nconst results = body.search(oldFragment, { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (const item of results.items) {\n  item.load(\"text\");\n  await context.sync();\n  item.insertText(item.text.split(oldFragment).join(newFragment), Word.InsertLocation.replace);\n  fixed = true;\n}\n\nawait context.sync();\n\nreturn { fixed };\n", "ps1": "# Fix a typo in the SmartArt (\"Archivist\") diagram embedded in the document:\n#   \"Aufsetzen des Systems auf projektfremden einem Rechner\"\n#     -> \"Aufsetzen des Systems auf projektfremdem Rechner\"\n#\n# The text lives inside the SmartArt graphic's diagram-data node, so the\n# canonical way to touch it through the Word object model is via\n# Shape.SmartArt.AllNodes(...).TextFrame2.TextRange.Text. We walk every\n# shape in the document, and for any shape that exposes a SmartArt\n# diagram, walk its nodes looking for the exact wording to repair.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Aufsetzen des Systems auf projektfremden einem Rechner\"\n$newText = \"Aufsetzen des Systems auf projektfremdem Rechner\"\n$oldFragment = \"projektfremden einem Rechner\"\n$newFragment = \"projektfremdem Rechner\"\n\nfunction Fix-NodeText($range) {\n    if ($range -eq $null) { return $false }\n    $t = $range.Text\n    if ([string]::IsNullOrEmpty($t)) { return $false }\n    if ($t -eq $oldText -or $t.Contains($oldFragment)) {\n        $range.Text = $t.Replace($oldFragment, $newFragment)\n        return $true\n    }\n    return $false\n}\n\n$fixed = $false\n\n# 1) Walk every shape on the page; if it carries a SmartArt diagram, walk\n#    every node in it and repair the wording wherever it is found.\nfor ($i = 1; $i -le $d.Shapes.Count; $i++) {\n    $shp = $d.Shapes.Item($i)\n\n    if ($shp.HasSmartArt) {\n        $sa = $shp.SmartArt\n        $nodes = $sa.AllNodes\n        for ($j = 1; $j -le $nodes.Count; $j++) {\n            $node = $nodes.Item($j)\n            if (Fix-NodeText $node.TextFrame2.TextRange) { $fixed = $true }\n        }\n    }\n\n    # Some hosts surface a SmartArt shape's cached text through the plain\n    # shape text frame as well - cover that path too.\n    if ($shp.TextFrame -ne $null -and $shp.TextFrame.HasText) {\n        if (Fix-NodeText $shp.TextFrame.TextRange) { $fixed = $true }\n    }\n}\n\n# 2) Belt-and-braces: also run an ordinary document-wide Find & Replace so\n#    the fix still lands if the phrase ever lives in normal story text.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldFragment\n$find.Replacement.Text = $newFragment\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $false\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\nif ($find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)) {\n    $fixed = $true\n}\n\nWrite-Output (\"fixed=\" + $fixed)\n"}
